$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.134157180786133
$ws.Range("B1").Value = 2.439523458480835
$ws.Range("C1").Value = 2.485842943191528
$ws.Range("D1").Value = 3.255719184875488
$ws.Range("E1").Value = 2.33658504486084
